$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / volume snapshot.
# Price-column (D) cells are stored as text in the source sheet, and some of the
# new values look exactly like plain numbers (trailing zeros such as "1.00", or
# multi-dot thousands grouping such as "66.054.71"). Excel's Range.Value setter
# auto-detects numeric-looking strings and would silently coerce them to Number,
# changing both the cell type and the displayed text (e.g. "10.90" -> 10.9).
# Force text storage via NumberFormat "@" before writing, then restore the
# default "Normal" style so no stray custom-format style is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.054.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.314.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.305.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.574"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.844.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "586.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.931.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.316.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.895"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.30%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "559.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.77%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.775.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0684"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.72%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.72%  "
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0408"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("E47").Value = "  -9.13%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.126"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.88%  "
